$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update res_line/pl_mw.xlsx style values for case with 380 kV (rows 2-25)
$rowData = @{
    2 = @(1.083068580487122, 0.3889912089761083, 0.0257865923873819, 0.07562253673345154, 1.260833122626792, 0.7608190033731219, 0.4319326158130963, 1.253478116961475)
    3 = @(0.9674094650524125, 0.3422684674269476, 0.02618341187955409, 0.07133779144054841, 1.217775299302218, 0.749140164128562, 0.3895921539997005, 1.269049954758387)
    4 = @(0.8967779099125437, 0.313687255838488, 0.02644378740342823, 0.06876489150391052, 1.192292039348601, 0.7425255544600091, 0.3638062761704433, 1.279164933352767)
    5 = @(0.8680899028225326, 0.3020659771148644, 0.02655406488512213, 0.06773075542664486, 1.182145330265087, 0.7399689340136888, 0.3533504694242993, 1.283425591205258)
    6 = @(0.863331970124932, 0.3001378004846913, 0.0265726273978828, 0.06755989741227708, 1.180474786978834, 0.7395527688929278, 0.3516174114867638, 1.284141433036609)
    7 = @(0.8963906310764855, 0.3135304241199606, 0.02644525779591334, 0.06875088701708165, 1.192154236666994, 0.7424905139467199, 0.3636650556602419, 1.279221833122946)
    8 = @(1.043108730473307, 0.3728584082916768, 0.02591992372500229, 0.07413298794923406, 1.245787428709249, 0.7566760899403562, 0.4172891464775219, 1.258732020899409)
    9 = @(1.333951536604218, 0.4901005027218162, 0.025023982524921, 0.08515814239237329, 1.358635580790988, 0.7889546083764714, 0.5241733850616583, 1.222969436414317)
    10 = @(1.549682480790921, 0.5768693767566333, 0.02444966655414049, 0.09356162687733161, 1.446373341710967, 0.8154579872626471, 0.6038330271779557, 1.19941986094782)
    11 = @(1.648301261758832, 0.6164983264724242, 0.02420708247706571, 0.09745405272481378, 1.487369310619016, 0.8281358075262375, 0.6403357010164825, 1.189305283742193)
    12 = @(1.685717092738514, 0.6315287491469803, 0.02411794732159578, 0.09893830691741812, 1.503051786114753, 0.8330271019375459, 0.6541977228741302, 1.185561801801185)
    13 = @(1.677655735945905, 0.6282906013385059, 0.0241370222685795, 0.09861818491866359, 1.499667210743382, 0.8319696318308729, 0.6512105245237194, 1.186364161444395)
    14 = @(1.651378050666608, 0.6177344046301982, 0.02419969440127545, 0.09757595572351363, 1.488656333238907, 0.8285363973645019, 0.6414753446965591, 1.188995564026492)
    15 = @(1.635291514350456, 0.6112715638602708, 0.02423843914073842, 0.09693890648997439, 1.481932524906483, 0.8264452591446769, 0.6355174115283404, 1.190618683370381)
    16 = @(1.543247354958282, 0.5742827978182277, 0.02446589924849718, 0.09330867140024424, 1.44371616569785, 0.8146420608913161, 0.6014529223013767, 1.20009296692205)
    17 = @(1.486905865424376, 0.5516325179771684, 0.02461025053670696, 0.09109966183652318, 1.420550883784529, 0.8075611152854378, 0.580624158679143, 1.206058825459166)
    18 = @(1.454544961506485, 0.5386194157021578, 0.02469503289430897, 0.08983563671677786, 1.407328678409044, 0.8035467535544853, 0.5686688887573155, 1.20954649875209)
    19 = @(1.443595830569336, 0.5342158990137023, 0.02472403899102105, 0.08940877461525076, 1.402869288389525, 0.8021975573637832, 0.5646252795236251, 1.210737009267518)
    20 = @(1.492898819256709, 0.5540421450086797, 0.02459470213551818, 0.0913341360459512, 1.423006309233557, 0.8083088405401071, 0.5828388325993359, 1.205417921215307)
    21 = @(1.65909450150906, 0.6208343601414299, 0.02418121176135202, 0.09788180267171498, 1.491886182556414, 0.8295423568430209, 0.6443337286200546, 1.18822029926023)
    22 = @(1.768127908675524, 0.664626213340739, 0.02392688152733058, 0.1022211083982185, 1.537826346321367, 0.8439476061414268, 0.6847535241899294, 1.17748626304499)
    23 = @(1.709896165592909, 0.6412405406789503, 0.02406115298308009, 0.09989955885980351, 1.513221954022299, 0.8362105694615991, 0.6631593567838365, 1.18316873594415)
    24 = @(1.490189309746768, 0.5529527251653121, 0.02460172598616417, 0.09122811169944356, 1.421895913097075, 0.8079706178408657, 0.581837517979892, 1.205707493972184)
    25 = @(1.254921411902217, 0.4582785970085865, 0.02525177308576154, 0.08212346921155955, 1.327270233735618, 0.7797377847152376, 0.4950653572952177, 1.232167932425583)
}

$cols = @("B","C","D","E","F","I","M","N")
foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

Write-Host "Updated $($rowData.Keys.Count) rows"